# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to the newly scraped totals.

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" (exhibitions) --
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 76
$ws1.Range("F3").Value = 3866
$ws1.Range("F4").Value = 2300
$ws1.Range("F5").Value = 455
$ws1.Range("F6").Value = 13
$ws1.Range("F10").Value = 105
$ws1.Range("F11").Value = 1431
$ws1.Range("F12").Value = 252
$ws1.Range("F13").Value = 2518
$ws1.Range("F14").Value = 178

# -- Sheet "全部类型" (all categories) --
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 76
$ws4.Range("F3").Value = 3866
$ws4.Range("F4").Value = 2300
$ws4.Range("F5").Value = 455
$ws4.Range("F6").Value = 13
$ws4.Range("F11").Value = 105
$ws4.Range("F14").Value = 1431
$ws4.Range("F15").Value = 252
$ws4.Range("F16").Value = 2518
$ws4.Range("F17").Value = 178
